$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Cells.Item(2, 1).Value = "ON"
$ws.Cells.Item(2, 2).Value = "SOFRRATE"
$ws.Cells.Item(2, 3).Value = "DEPOSIT"
$ws.Cells.Item(2, 4).Value = 0.0436

$ws.Cells.Item(3, 1).Value = "3M"
$ws.Cells.Item(3, 2).Value = "SQZ25"
$ws.Cells.Item(3, 3).Value = "FUTURE"
$ws.Cells.Item(3, 4).Value = 96.23

$ws.Cells.Item(4, 1).Value = "4M"
$ws.Cells.Item(4, 2).Value = "SQF26"
$ws.Cells.Item(4, 3).Value = "FUTURE"
$ws.Cells.Item(4, 4).Value = 96.335

$ws.Cells.Item(5, 1).Value = "5M"
$ws.Cells.Item(5, 2).Value = "SQG26"
$ws.Cells.Item(5, 3).Value = "FUTURE"
$ws.Cells.Item(5, 4).Value = 96.415

$ws.Cells.Item(6, 1).Value = "6M"
$ws.Cells.Item(6, 2).Value = "SQH26"
$ws.Cells.Item(6, 3).Value = "FUTURE"
$ws.Cells.Item(6, 4).Value = 96.49

$ws.Cells.Item(7, 1).Value = "9M"
$ws.Cells.Item(7, 2).Value = "SQM26"
$ws.Cells.Item(7, 3).Value = "FUTURE"
$ws.Cells.Item(7, 4).Value = 96.735

$ws.Cells.Item(8, 1).Value = "12M"
$ws.Cells.Item(8, 2).Value = "SQU26"
$ws.Cells.Item(8, 3).Value = "FUTURE"
$ws.Cells.Item(8, 4).Value = 96.895

$ws.Cells.Item(9, 1).Value = "15M"
$ws.Cells.Item(9, 2).Value = "SQZ26"
$ws.Cells.Item(9, 3).Value = "FUTURE"
$ws.Cells.Item(9, 4).Value = 96.975

$ws.Cells.Item(10, 1).Value = "0M"
$ws.Cells.Item(10, 2).Value = "SQU25"
$ws.Cells.Item(10, 3).Value = "FUTURE"
$ws.Cells.Item(10, 4).Value = 95.895

$ws.Cells.Item(11, 1).Value = "1M"
$ws.Cells.Item(11, 2).Value = "SQV25"
$ws.Cells.Item(11, 3).Value = "FUTURE"
$ws.Cells.Item(11, 4).Value = 96.005

$ws.Cells.Item(12, 1).Value = "2M"
$ws.Cells.Item(12, 2).Value = "SQX25"
$ws.Cells.Item(12, 3).Value = "FUTURE"
$ws.Cells.Item(12, 4).Value = 96.125

$ws.Cells.Item(13, 1).Value = "2Y"
$ws.Cells.Item(13, 2).Value = "SOFROIS"
$ws.Cells.Item(13, 3).Value = "OIS"
$ws.Cells.Item(13, 4).Value = 0.034355

$ws.Cells.Item(14, 1).Value = "3Y"
$ws.Cells.Item(14, 2).Value = "SOFROIS"
$ws.Cells.Item(14, 3).Value = "OIS"
$ws.Cells.Item(14, 4).Value = 0.033474

$ws.Cells.Item(15, 1).Value = "5Y"
$ws.Cells.Item(15, 2).Value = "SOFROIS"
$ws.Cells.Item(15, 3).Value = "OIS"
$ws.Cells.Item(15, 4).Value = 0.034039

$ws.Cells.Item(16, 1).Value = "10Y"
$ws.Cells.Item(16, 2).Value = "SOFROIS"
$ws.Cells.Item(16, 3).Value = "OIS"
$ws.Cells.Item(16, 4).Value = 0.037427

$ws.Cells.Item(17, 1).Value = "20Y"
$ws.Cells.Item(17, 2).Value = "SOFROIS"
$ws.Cells.Item(17, 3).Value = "OIS"
$ws.Cells.Item(17, 4).Value = 0.04088

$ws.Cells.Item(18, 1).Value = "30Y"
$ws.Cells.Item(18, 2).Value = "SOFROIS"
$ws.Cells.Item(18, 3).Value = "OIS"
$ws.Cells.Item(18, 4).Value = 0.04067
